$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.938.49"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.576.87"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "302.77"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "94.56"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("E13").Value = "  +6.47%  "
$ws.Range("D14").Value = "2.548.56"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "0.886"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "14.27"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "42.998.73"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "12.99"
$ws.Range("E18").Value = "  +5.46%  "
$ws.Range("D19").Value = "0.0₃0997"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "72.01"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "254.03"
$ws.Range("E22").Value = "  -4.34%  "
$ws.Range("D23").Value = "2.96"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  -3.71%  "
$ws.Range("D25").Value = "28.87"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "10.33"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "37.51"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -5.49%  "
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "154.94"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "3.42"
$ws.Range("E33").Value = "  -5.45%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "0.0804"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").Value = "18.34"
$ws.Range("E36").Value = "  +9.96%  "
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").Value = "23.38"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").Value = "2.22"
$ws.Range("E40").Value = "  +36.31%  "
$ws.Range("D41").Value = "3.44"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "2.082.46"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "9.28"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").Value = "85.42"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").Value = "76.61"
$ws.Range("E48").Value = "  +10.77%  "
$ws.Range("D49").Value = "106.76"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "2.818.16"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").Value = "1.69"
$ws.Range("E51").Value = "  +2.54%  "
